# Fixes for View Payments
# Adds 5 new SQL rows (144-148) to the "SQL" sheet describing "last 90 days"
# style queries (archive/active/NPI null/not-null variants).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL")

# --- Build the 4 distinct SQL query strings (CRLF line endings, matching
# the convention already used throughout this workbook's shared strings) ---

$qNpiNotNull = "Select p.prov_tax_id_nbr" + "`r`n" + `
  "from PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p, OLE.PROC_CTL PC WHERE PC.PROC_CTL_ID = CP.PROC_CTL_ID" + "`r`n" + `
  "and cp.prov_key_id = p.prov_key_id AND PC.EXTRACT_STS_CD = 'C'" + "`r`n" + `
  "and cp.setl_dt > current date - 90 days and p.PROV_NPI_NBR is not null" + "`r`n" + `
  "group by p.prov_tax_id_nbr" + "`r`n" + `
  "having count(*) between 1 and 30" + "`r`n" + `
  "order by count(*) desc" + "`r`n" + `
  "fetch first row only"

$qNpiNull = "Select p.prov_tax_id_nbr" + "`r`n" + `
  "from PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p, OLE.PROC_CTL PC WHERE PC.PROC_CTL_ID = CP.PROC_CTL_ID" + "`r`n" + `
  "and cp.prov_key_id = p.prov_key_id AND PC.EXTRACT_STS_CD = 'C'" + "`r`n" + `
  "and cp.setl_dt > current date - 90 days and p.PROV_NPI_NBR is null" + "`r`n" + `
  "group by p.prov_tax_id_nbr" + "`r`n" + `
  "having count(*) between 1 and 30" + "`r`n" + `
  "order by count(*) desc" + "`r`n" + `
  "fetch first row only with ur"

$qArchvY = "Select p.prov_tax_id_nbr" + "`r`n" + `
  "from PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p, OLE.PROC_CTL PC WHERE PC.PROC_CTL_ID = CP.PROC_CTL_ID" + "`r`n" + `
  "and cp.prov_key_id = p.prov_key_id AND PC.EXTRACT_STS_CD = 'C'" + "`r`n" + `
  "and cp.setl_dt > current date - 90 days and cp.ARCHV_IND='Y'" + "`r`n" + `
  "group by p.prov_tax_id_nbr" + "`r`n" + `
  "having count(*) between 1 and 30" + "`r`n" + `
  "order by count(*) desc" + "`r`n" + `
  "fetch first row only with ur"

$qArchvNotIn = "Select p.prov_tax_id_nbr" + "`r`n" + `
  "from PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p, OLE.PROC_CTL PC WHERE PC.PROC_CTL_ID = CP.PROC_CTL_ID" + "`r`n" + `
  "and cp.prov_key_id = p.prov_key_id AND PC.EXTRACT_STS_CD = 'C'" + "`r`n" + `
  "and cp.setl_dt > current date - 90 days and cp.ARCHV_IND not in ('Y')" + "`r`n" + `
  "group by p.prov_tax_id_nbr" + "`r`n" + `
  "having count(*) between 1 and 30" + "`r`n" + `
  "order by count(*) desc" + "`r`n" + `
  "fetch first row only with ur"

# --- New row data: Sno (A), Query (B), Comments (C) ---
# Row 144: Sno 143, npi-not-null query, "Last 90 days"
# Row 145: Sno 144, npi-null query, "Npi last 90 days"
# Row 146: Sno 145, archive='Y' query, "archive only 90 days"
# Row 147: Sno 146, archive not-in query, "active only 90 days"
# Row 148: Sno 147, npi-not-null query (reused), "Npi last 90 days" (reused)

$rows = 144, 145, 146, 147, 148
$snoValues = "143", "144", "145", "146", "147"
$queryValues = $qNpiNotNull, $qNpiNull, $qArchvY, $qArchvNotIn, $qNpiNotNull
$commentValues = "Last 90 days", "Npi last 90 days", "archive only 90 days", "active only 90 days", "Npi last 90 days"

$srcA = $ws.Range("A143")
$srcB = $ws.Range("B143")
$srcC = $ws.Range("C143")

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]

    $cellA = $ws.Range("A$r")
    $cellA.Value = "'" + $snoValues[$i]
    $srcA.Copy()
    $cellA.PasteSpecial(-4122)

    $cellB = $ws.Range("B$r")
    $cellB.Value = $queryValues[$i]
    $srcB.Copy()
    $cellB.PasteSpecial(-4122)

    $cellC = $ws.Range("C$r")
    $cellC.Value = $commentValues[$i]
    $srcC.Copy()
    $cellC.PasteSpecial(-4122)

    $ws.Rows.Item($r).RowHeight = 129.6
}

# Move selection/view to the newly-added last row, matching the authored
# sheetView state (topLeftCell scrolled down, C148 selected).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 146
$win.ScrollColumn = 1
[void]$ws.Range("C148").Select()

Write-Host "Added rows 144-148 to SQL sheet"
